{"js": "// Insert \"  \u0631\u0641\u06cc\u0639\u06cc \" (two leading spaces, \"\u0631\u0641\u06cc\u0639\u06cc\", one trailing space) as\n// three new runs at the very start of the first paragraph, immediately\n// before the existing \"\u062d\u0645\u06cc\u062f\u0631\u0636\u0627\" run \u2014 matching the authored diff exactly,\n// including per-run formatting (the middle run is RTL with a CS font hint,\n// the space runs are plain/LTR).\n//\n// We build the new runs via a Flat-OPC `insertOoxml` payload so that the\n// exact run-level properties (w:rtl, w:hint=\"cs\", w:lang) land precisely as\n// authored, rather than relying on font-property heuristics that Office.js's\n// higher level `font.*` setters do not fully expose (there is no direct\n// \"rtl\"/\"hint\" setter on Range.font in the Word JS API surface).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Collapsed range at the very beginning of the paragraph (right after\n// <w:pPr>, right before the existing run).\nconst insertionPoint = firstParagraph.getRange(\"Start\");\n\nconst flatOpcPayload =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<?mso-application progid=\"Word.Document\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n      '<pkg:xmlData>' +\n        '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n          '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n        '</Relationships>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:r>' +\n                '<w:rPr>' +\n                  '<w:rFonts w:ascii=\"IRANYekan Medium\" w:hAnsi=\"IRANYekan Medium\" w:cs=\"IRANYekan Medium\"/>' +\n                  '<w:sz w:val=\"28\"/>' +\n                  '<w:szCs w:val=\"28\"/>' +\n                  '<w:lang w:bidi=\"fa-IR\"/>' +\n                '</w:rPr>' +\n                '<w:t xml:space=\"preserve\">  </w:t>' +\n              '</w:r>' +\n              '<w:r>' +\n                '<w:rPr>' +\n                  '<w:rFonts w:ascii=\"IRANYekan Medium\" w:hAnsi=\"IRANYekan Medium\" w:cs=\"IRANYekan Medium\" w:hint=\"cs\"/>' +\n                  '<w:sz w:val=\"28\"/>' +\n                  '<w:szCs w:val=\"28\"/>' +\n                  '<w:rtl/>' +\n                  '<w:lang w:bidi=\"fa-IR\"/>' +\n                '</w:rPr>' +\n                '<w:t>\\u0631\\u0641\\u06CC\\u0639\\u06CC</w:t>' +\n              '</w:r>' +\n              '<w:r>' +\n                '<w:rPr>' +\n                  '<w:rFonts w:ascii=\"IRANYekan Medium\" w:hAnsi=\"IRANYekan Medium\" w:cs=\"IRANYekan Medium\"/>' +\n                  '<w:sz w:val=\"28\"/>' +\n                  '<w:szCs w:val=\"28\"/>' +\n                  '<w:lang w:bidi=\"fa-IR\"/>' +\n                '</w:rPr>' +\n                '<w:t xml:space=\"preserve\"> </w:t>' +\n              '</w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\ninsertionPoint.insertOoxml(flatOpcPayload, \"Start\");\nawait context.sync();\n", "ps1": "# Insert \"  \u0631\u0641\u06cc\u0639\u06cc \" (two leading spaces, \"\u0631\u0641\u06cc\u0639\u06cc\", one trailing space) as\n# three new runs at the very start of the first paragraph, immediately\n# before the existing \"\u062d\u0645\u06cc\u062f\u0631\u0636\u0627\" run \u2014 matching the authored diff exactly,\n# including per-run formatting (the middle run is RTL with a CS font hint,\n# the space runs are plain/LTR).\n#\n# Word COM's Range object does not expose a simple per-run \"hint\"/\"rtl\"\n# property setter either, so \u2014 exactly like the Office.js version of this\n# edit \u2014 we build the new runs as a small Flat-OPC WordProcessingML package\n# and hand it to Range.InsertXML, which lets us land the exact authored\n# run-level XML (w:rtl, w:hint=\"cs\", w:lang) instead of reconstructing it\n# through higher-level formatting properties.\n\n$d = $word.ActiveDocument\n$p = $d.Paragraphs(1)\n$r = $p.Range\n$r.SetRange($r.Start, $r.Start)\n\n$flatOpc = '<?xml version=\"1.0\" standalone=\"yes\"?><?mso-application progid=\"Word.Document\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii=\"IRANYekan Medium\" w:hAnsi=\"IRANYekan Medium\" w:cs=\"IRANYekan Medium\"/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/><w:lang w:bidi=\"fa-IR\"/></w:rPr><w:t xml:space=\"preserve\">  </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"IRANYekan Medium\" w:hAnsi=\"IRANYekan Medium\" w:cs=\"IRANYekan Medium\" w:hint=\"cs\"/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/><w:rtl/><w:lang w:bidi=\"fa-IR\"/></w:rPr><w:t>\u0631\u0641\u06cc\u0639\u06cc</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"IRANYekan Medium\" w:hAnsi=\"IRANYekan Medium\" w:cs=\"IRANYekan Medium\"/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/><w:lang w:bidi=\"fa-IR\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$r.InsertXML($flatOpc)\n"}
